$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.624.03'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '1.862.61'
$ws.Range("E3").Value = '  -0.99%  '
$ws.Range("E4").Value = '  +0.43%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '334.21'
$ws.Range("E5").Value = '  +0.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.011'
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4686'
$ws.Range("E7").Value = '  -0.46%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3908'
$ws.Range("E8").Value = '  -0.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.67'
$ws.Range("E9").Value = '  -4.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07992'
$ws.Range("E10").Value = '  -0.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.002'
$ws.Range("E11").Value = '  -2.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.74'
$ws.Range("E12").Value = '  -1.99%  '
$ws.Range("D13").Value = '1.874.36'
$ws.Range("E13").Value = '  -0.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.989'
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("E15").Value = '  +1.53%  '
$ws.Range("E16").Value = '  +0.38%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.16'
$ws.Range("E17").Value = '  +1.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06708'
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("E19").Value = '  -0.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.06'
$ws.Range("E20").Value = '  -1.61%  '
$ws.Range("E21").Value = '  +0.48%  '
$ws.Range("D22").Value = '27.607.96'
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.451'
$ws.Range("E23").Value = '  -1.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.88'
$ws.Range("E24").Value = '  -1.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.308'
$ws.Range("E25").Value = '  -0.18%  '
$ws.Range("D26").Value = '2.097.92'
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '159.73'
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.81'
$ws.Range("E28").Value = '  -2.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.135'
$ws.Range("E29").Value = '  +1.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.425'
$ws.Range("E30").Value = '  -2.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '121.47'
$ws.Range("E31").Value = '  -0.41%  '
$ws.Range("E32").Value = '  -0.75%  '
$ws.Range("E33").Value = '  -0.29%  '
$ws.Range("E34").Value = '  +0.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.291'
$ws.Range("E35").Value = '  -1.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.332'
$ws.Range("E36").Value = '  -8.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06053'
$ws.Range("E37").Value = '  -1.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02231'
$ws.Range("E38").Value = '  -1.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.343'
$ws.Range("E39").Value = '  +2.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.192'
$ws.Range("E40").Value = '  -2.98%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.010'
$ws.Range("E41").Value = '  +0.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5942'
$ws.Range("E42").Value = '  -0.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1883'
$ws.Range("E43").Value = '  -0.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.24'
$ws.Range("E44").Value = '  -0.79%  '
$ws.Range("E45").Value = '  -0.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5633'
$ws.Range("E46").Value = '  -1.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.22'
$ws.Range("E47").Value = '  +0.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.921'
$ws.Range("E48").Value = '  -1.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.206'
$ws.Range("E49").Value = '  -5.57%  '
$ws.Range("E50").Value = '  -2.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '112.42'
$ws.Range("E51").Value = '  -1.97%  '
